$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 3125.6428
$ws.Range("I33").Value = 3913.9092
$ws.Range("J33").Value = 235.33333
$ws.Range("K33").Value = 3913.9092
$ws.Range("L33").Value = 235.33333
$ws.Range("M33").Value = -3684.9092
$ws.Range("N33").Value = -693.3333299999999
$ws.Range("H64").Value = 3000
$ws.Range("J64").Value = 0
$ws.Range("L64").Value = 0
$ws.Range("N64").ClearContents()
$ws.Range("H67").Value = 3000
$ws.Range("J67").Value = 0
$ws.Range("L67").Value = 0
$ws.Range("N67").ClearContents()
$ws.Range("H74").Value = 7800
$ws.Range("I74").Value = 7800
$ws.Range("J74").Value = 7800
$ws.Range("K74").Value = 7800
$ws.Range("L74").Value = 7800
$ws.Range("M74").Value = -6864
$ws.Range("N74").Value = -9672
$ws.Range("H77").Value = 7800
$ws.Range("I77").Value = 7800
$ws.Range("J77").Value = 7800
$ws.Range("K77").Value = 39000
$ws.Range("L77").Value = 39000
$ws.Range("M77").Value = -34320
$ws.Range("N77").Value = -48360
$ws.Range("H98").Value = 25430.295
$ws.Range("I98").Value = 23488.133
$ws.Range("K98").Value = 23488.133
$ws.Range("M98").Value = -21990.133
$ws.Range("H122").Value = 25430.295
$ws.Range("I122").Value = 23488.133
$ws.Range("K122").Value = 70464.399
$ws.Range("M122").Value = -68014.399
$ws.Range("H132").Value = 3669.875
$ws.Range("I132").Value = 3740.7886
$ws.Range("K132").Value = 11222.3658
$ws.Range("M132").Value = -8692.3658
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4114.3096
$ws.Range("I32").Value = 2230.077
$ws.Range("J32").Value = 7176.1875
$ws.Range("K32").Value = 2230.077
$ws.Range("L32").Value = 7176.1875
$ws.Range("M32").Value = -1943.077
$ws.Range("N32").Value = -7750.1875
$ws.Range("H48").Value = 247950
$ws.Range("J48").Value = 247950
$ws.Range("L48").Value = 247950
$ws.Range("N48").Value = -248718
$ws.Range("H74").Value = 44957.387
$ws.Range("J74").Value = 233818.75
$ws.Range("L74").Value = 233818.75
$ws.Range("N74").Value = -235566.75
$ws.Range("H77").Value = 44957.387
$ws.Range("J77").Value = 233818.75
$ws.Range("L77").Value = 1169093.75
$ws.Range("N77").Value = -1177829.75
$ws.Range("H97").Value = 10157.708
$ws.Range("I97").Value = 7355.278
$ws.Range("K97").Value = 7355.278
$ws.Range("M97").Value = -6859.278
$ws.Range("H110").Value = 1431
$ws.Range("I110").Value = 917.2
$ws.Range("K110").Value = 917.2
$ws.Range("M110").Value = 1127.8
$ws.Range("H122").Value = 3112.5
$ws.Range("I122").Value = 2177.3572
$ws.Range("K122").Value = 6532.071599999999
$ws.Range("M122").Value = -4082.071599999999
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 999.9
$ws.Range("I22").Value = 888.7778
$ws.Range("K22").Value = 888.7778
$ws.Range("M22").Value = -715.7778
$ws.Range("H56").Value = 18666.666
$ws.Range("J56").Value = 18666.666
$ws.Range("L56").Value = 18666.666
$ws.Range("N56").Value = -20144.666
$ws.Range("H134").Value = 7665.6787
$ws.Range("I134").Value = 7312.5557
$ws.Range("J134").Value = 8301.299999999999
$ws.Range("K134").Value = 21937.6671
$ws.Range("L134").Value = 24903.9
$ws.Range("M134").Value = -19402.6671
$ws.Range("N134").Value = -29973.9
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 249.12
$ws.Range("I7").Value = 140.05882
$ws.Range("K7").Value = 140.05882
$ws.Range("M7").Value = -27.05882
$ws.Range("H22").Value = 1142.7273
$ws.Range("I22").Value = 703.3570999999999
$ws.Range("J22").Value = 1911.625
$ws.Range("K22").Value = 703.3570999999999
$ws.Range("L22").Value = 1911.625
$ws.Range("M22").Value = -353.3570999999999
$ws.Range("N22").Value = -2611.625
$ws.Range("H31").Value = 26011.979
$ws.Range("I31").Value = 12237.7
$ws.Range("J31").Value = 29947.486
$ws.Range("K31").Value = 12237.7
$ws.Range("L31").Value = 29947.486
$ws.Range("M31").Value = -11942.7
$ws.Range("N31").Value = -30537.486
$ws.Range("H34").Value = 26011.979
$ws.Range("I34").Value = 12237.7
$ws.Range("J34").Value = 29947.486
$ws.Range("K34").Value = 12237.7
$ws.Range("L34").Value = 29947.486
$ws.Range("M34").Value = -12035.7
$ws.Range("N34").Value = -30351.486
$ws.Range("H58").Value = 3715.5715
$ws.Range("J58").Value = 4104.6
$ws.Range("L58").Value = 4104.6
$ws.Range("N58").Value = -4510.6
$ws.Range("H62").Value = 2500
$ws.Range("I62").Value = 1500
$ws.Range("K62").Value = 1500
$ws.Range("M62").Value = -876
$ws.Range("H65").Value = 2500
$ws.Range("I65").Value = 1500
$ws.Range("K65").Value = 7500
$ws.Range("M65").Value = -4380
$ws.Range("H86").Value = 14019.042
$ws.Range("I86").Value = 13779.111
$ws.Range("K86").Value = 13779.111
$ws.Range("M86").Value = -12656.111
$ws.Range("H89").Value = 14019.042
$ws.Range("I89").Value = 13779.111
$ws.Range("K89").Value = 68895.55500000001
$ws.Range("M89").Value = -63279.55500000001
$ws.Range("H107").Value = 488.44446
$ws.Range("I107").Value = 286.33334
$ws.Range("K107").Value = 286.33334
$ws.Range("M107").Value = 1633.66666
$ws.Range("H122").Value = 3346.1428
$ws.Range("I122").Value = 2685.6
$ws.Range("K122").Value = 8056.799999999999
$ws.Range("M122").Value = -5606.799999999999
$ws.Range("H136").Value = 3715.5715
$ws.Range("J136").Value = 4104.6
$ws.Range("L136").Value = 12313.8
$ws.Range("N136").Value = -17413.8
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 6193580.5
$ws.Range("J4").Value = 359903.1
$ws.Range("L4").Value = 1079709.3
$ws.Range("N4").Value = -1079933.3
$ws.Range("H23").Value = 993.3333
$ws.Range("J23").Value = 993.3333
$ws.Range("L23").Value = 2979.9999
$ws.Range("N23").Value = -3449.9999
$ws.Range("H34").Value = 2200.4
$ws.Range("J34").Value = 3109.7
$ws.Range("L34").Value = 9329.099999999999
$ws.Range("N34").Value = -9497.099999999999
$ws.Range("H132").Value = 2057
$ws.Range("I132").Value = 1649.25
$ws.Range("J132").Value = 2600.6667
$ws.Range("K132").Value = 14843.25
$ws.Range("L132").Value = 23406.0003
$ws.Range("M132").Value = -12313.25
$ws.Range("N132").Value = -28466.0003
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H4").Value = 4401.5
$ws.Range("J4").Value = 0
$ws.Range("L4").Value = 0
$ws.Range("N4").ClearContents()
$ws.Range("H44").Value = 10500
$ws.Range("I44").Value = 10500
$ws.Range("K44").Value = 10500
$ws.Range("M44").Value = -9904
$ws.Range("H92").Value = 2500000
$ws.Range("J92").Value = 2500000
$ws.Range("L92").Value = 2500000
$ws.Range("N92").Value = -2503744
$ws.Range("H97").Value = 2085.9
$ws.Range("I97").Value = 1718.64
$ws.Range("J97").Value = 3922.2
$ws.Range("K97").Value = 1718.64
$ws.Range("L97").Value = 3922.2
$ws.Range("M97").Value = -1222.64
$ws.Range("N97").Value = -4914.2
$ws.Range("H101").Value = 30000.5
$ws.Range("J101").Value = 30000.5
$ws.Range("L101").Value = 30000.5
$ws.Range("N101").Value = -36490.5
$ws.Range("H102").Value = 34956.547
$ws.Range("I102").Value = 2063.45
$ws.Range("K102").Value = 2063.45
$ws.Range("M102").Value = -441.4499999999998
$ws.Range("H126").Value = 3777.1924
$ws.Range("J126").Value = 4642.778
$ws.Range("L126").Value = 13928.334
$ws.Range("N126").Value = -18868.334
$ws.Range("H132").Value = 5231.174
$ws.Range("I132").Value = 3133.3489
$ws.Range("J132").Value = 35300
$ws.Range("K132").Value = 9400.046699999999
$ws.Range("L132").Value = 105900
$ws.Range("M132").Value = -6870.046699999999
$ws.Range("N132").Value = -110960
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1033
$ws.Range("I16").Value = 870.2
$ws.Range("K16").Value = 870.2
$ws.Range("M16").Value = -700.2
$ws.Range("H22").Value = 75125.336
$ws.Range("I22").Value = 128012.14
$ws.Range("J22").Value = 1083.8
$ws.Range("K22").Value = 128012.14
$ws.Range("L22").Value = 1083.8
$ws.Range("M22").Value = -127717.14
$ws.Range("N22").Value = -1673.8
$ws.Range("H27").Value = 75125.336
$ws.Range("I27").Value = 128012.14
$ws.Range("J27").Value = 1083.8
$ws.Range("K27").Value = 128012.14
$ws.Range("L27").Value = 1083.8
$ws.Range("M27").Value = -127905.14
$ws.Range("N27").Value = -1297.8
$ws.Range("H46").Value = 10064.471
$ws.Range("I46").Value = 8357
$ws.Range("J46").Value = 11259.7
$ws.Range("K46").Value = 8357
$ws.Range("L46").Value = 11259.7
$ws.Range("M46").Value = -8169
$ws.Range("N46").Value = -11635.7
$ws.Range("H55").Value = 1047.5676
$ws.Range("I55").Value = 633.4761999999999
$ws.Range("J55").Value = 1591.0625
$ws.Range("K55").Value = 633.4761999999999
$ws.Range("L55").Value = 1591.0625
$ws.Range("M55").Value = -460.4761999999999
$ws.Range("N55").Value = -1937.0625
$ws.Range("H82").Value = 1115.6316
$ws.Range("I82").Value = 1116.5
$ws.Range("J82").Value = 1100
$ws.Range("K82").Value = 1116.5
$ws.Range("L82").Value = 1100
$ws.Range("M82").Value = -755.5
$ws.Range("N82").Value = -1822
$ws.Range("H85").Value = 1115.6316
$ws.Range("I85").Value = 1116.5
$ws.Range("J85").Value = 1100
$ws.Range("K85").Value = 1116.5
$ws.Range("L85").Value = 1100
$ws.Range("M85").Value = 131.5
$ws.Range("N85").Value = -3596
$ws.Range("H93").Value = 4463.857
$ws.Range("I93").Value = 4033.3333
$ws.Range("J93").Value = 4786.75
$ws.Range("K93").Value = 4033.3333
$ws.Range("L93").Value = 4786.75
$ws.Range("M93").Value = -2785.3333
$ws.Range("N93").Value = -7282.75
$ws.Range("H131").Value = 130000
$ws.Range("J131").Value = 130000
$ws.Range("L131").Value = 130000
$ws.Range("N131").Value = -140080
$ws.Range("H132").Value = 5000
$ws.Range("I132").Value = 5000
$ws.Range("K132").Value = 15000
$ws.Range("M132").Value = -12470
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H31").Value = 16499.5
$ws.Range("H62").Value = 15618.546
$ws.Range("J62").Value = 7537.875
$ws.Range("L62").Value = 7537.875
$ws.Range("N62").Value = -8785.875
$ws.Range("H65").Value = 15618.546
$ws.Range("J65").Value = 7537.875
$ws.Range("L65").Value = 37689.375
$ws.Range("N65").Value = -43929.375
$ws.Range("H69").Value = 23666
$ws.Range("J69").Value = 23666
$ws.Range("L69").Value = 23666
$ws.Range("N69").Value = -25164
$ws.Range("H72").Value = 23666
$ws.Range("J72").Value = 23666
$ws.Range("L72").Value = 70998
$ws.Range("N72").Value = -78486
$ws.Range("H81").Value = 1881.1666
$ws.Range("I81").Value = 1599.6666
$ws.Range("J81").Value = 2162.6667
$ws.Range("K81").Value = 3199.3332
$ws.Range("L81").Value = 4325.3334
$ws.Range("M81").Value = -2138.3332
$ws.Range("N81").Value = -6447.3334
$ws.Range("H84").Value = 1881.1666
$ws.Range("I84").Value = 1599.6666
$ws.Range("J84").Value = 2162.6667
$ws.Range("K84").Value = 15996.666
$ws.Range("L84").Value = 21626.667
$ws.Range("M84").Value = -10692.666
$ws.Range("N84").Value = -32234.667
$ws.Range("H100").Value = 1156
$ws.Range("I100").Value = 686.6667
$ws.Range("J100").Value = 1357.1428
$ws.Range("K100").Value = 1373.3334
$ws.Range("L100").Value = 2714.2856
$ws.Range("M100").Value = -832.3334
$ws.Range("N100").Value = -3796.2856
$ws.Range("H132").Value = 1850327.5
$ws.Range("I132").Value = 53740
$ws.Range("J132").Value = 2748621.2
$ws.Range("K132").Value = 161220
$ws.Range("L132").Value = 8245863.600000001
$ws.Range("M132").Value = -158690
$ws.Range("N132").Value = -8250923.600000001
$ws.Range("H136").Value = 5252.452
$ws.Range("I136").Value = 6998.0386
$ws.Range("K136").Value = 20994.1158
$ws.Range("M136").Value = -18444.1158
